$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

# Header row (text content unchanged; kept explicit for clarity)
$ws.Range("A1").Value = 'date'
$ws.Range("B1").Value = 'albania_gasoline'
$ws.Range("C1").Value = 'albania_diesel'
$ws.Range("D1").Value = 'albania_gas'

# Data rows: rolling 31-day window advanced forward (new prices pulled in)
# Row 2
$ws.Range("A2").Value = '26.01.2025'
$ws.Range("B2").Value = '3,44 '
$ws.Range("C2").Value = '3,44 '
$ws.Range("D2").Value = '1,13 '

# Row 3
$ws.Range("A3").Value = '27.01.2025'
$ws.Range("B3").Value = '3,43 '
$ws.Range("C3").Value = '3,43 '
$ws.Range("D3").Value = '1,13 '

# Row 4
$ws.Range("A4").Value = '28.01.2025'
$ws.Range("B4").Value = '3,42 '
$ws.Range("C4").Value = '3,42 '
$ws.Range("D4").Value = '1,13 '

# Row 5
$ws.Range("A5").Value = '29.01.2025'
$ws.Range("B5").Value = '3,43 '
$ws.Range("C5").Value = '3,43 '
$ws.Range("D5").Value = '1,13 '

# Row 6
$ws.Range("A6").Value = '30.01.2025'
$ws.Range("B6").Value = '3,43 '
$ws.Range("C6").Value = '3,43 '
$ws.Range("D6").Value = '1,13 '

# Row 7
$ws.Range("A7").Value = '31.01.2025'
$ws.Range("B7").Value = '3,43 '
$ws.Range("C7").Value = '3,43 '
$ws.Range("D7").Value = '1,13 '

# Row 8
$ws.Range("A8").Value = '01.02.2025'
$ws.Range("B8").Value = '3,44 '
$ws.Range("C8").Value = '3,44 '
$ws.Range("D8").Value = '1,13 '

# Row 9
$ws.Range("A9").Value = '02.02.2025'
$ws.Range("B9").Value = '3,44 '
$ws.Range("C9").Value = '3,44 '
$ws.Range("D9").Value = '1,13 '

# Row 10
$ws.Range("A10").Value = '03.02.2025'
$ws.Range("B10").Value = '3,41 '
$ws.Range("C10").Value = '3,41 '
$ws.Range("D10").Value = '1,12 '

# Row 11
$ws.Range("A11").Value = '04.02.2025'
$ws.Range("B11").Value = '3,43 '
$ws.Range("C11").Value = '3,43 '
$ws.Range("D11").Value = '1,13 '

# Row 12
$ws.Range("A12").Value = '05.02.2025'
$ws.Range("B12").Value = '3,40 '
$ws.Range("C12").Value = '3,40 '
$ws.Range("D12").Value = '1,12 '

# Row 13
$ws.Range("A13").Value = '06.02.2025'
$ws.Range("B13").Value = '3,39 '
$ws.Range("C13").Value = '3,39 '
$ws.Range("D13").Value = '1,12 '

# Row 14
$ws.Range("A14").Value = '07.02.2025'
$ws.Range("B14").Value = '3,40 '
$ws.Range("C14").Value = '3,40 '
$ws.Range("D14").Value = '1,12 '

# Row 15
$ws.Range("A15").Value = '08.02.2025'
$ws.Range("B15").Value = '3,40 '
$ws.Range("C15").Value = '3,40 '
$ws.Range("D15").Value = '1,12 '

# Row 16
$ws.Range("A16").Value = '09.02.2025'
$ws.Range("B16").Value = '3,40 '
$ws.Range("C16").Value = '3,40 '
$ws.Range("D16").Value = '1,12 '

# Row 17
$ws.Range("A17").Value = '10.02.2025'
$ws.Range("B17").Value = '3,41 '
$ws.Range("C17").Value = '3,41 '
$ws.Range("D17").Value = '1,12 '

# Row 18
$ws.Range("A18").Value = '11.02.2025'
$ws.Range("B18").Value = '3,42 '
$ws.Range("C18").Value = '3,42 '
$ws.Range("D18").Value = '1,13 '

# Row 19
$ws.Range("A19").Value = '12.02.2025'
$ws.Range("B19").Value = '3,42 '
$ws.Range("C19").Value = '3,42 '
$ws.Range("D19").Value = '1,13 '

# Row 20
$ws.Range("A20").Value = '13.02.2025'
$ws.Range("B20").Value = '3,41 '
$ws.Range("C20").Value = '3,41 '
$ws.Range("D20").Value = '1,12 '

# Row 21
$ws.Range("A21").Value = '14.02.2025'
$ws.Range("B21").Value = '3,40 '
$ws.Range("C21").Value = '3,40 '
$ws.Range("D21").Value = '1,12 '

# Row 22
$ws.Range("A22").Value = '15.02.2025'
$ws.Range("B22").Value = '3,40 '
$ws.Range("C22").Value = '3,40 '
$ws.Range("D22").Value = '1,12 '

# Row 23
$ws.Range("A23").Value = '16.02.2025'
$ws.Range("B23").Value = '3,40 '
$ws.Range("C23").Value = '3,40 '
$ws.Range("D23").Value = '1,12 '

# Row 24
$ws.Range("A24").Value = '17.02.2025'
$ws.Range("B24").Value = '3,39 '
$ws.Range("C24").Value = '3,39 '
$ws.Range("D24").Value = '1,12 '

# Row 25
$ws.Range("A25").Value = '18.02.2025'
$ws.Range("B25").Value = '3,39 '
$ws.Range("C25").Value = '3,39 '
$ws.Range("D25").Value = '1,12 '

# Row 26
$ws.Range("A26").Value = '19.02.2025'
$ws.Range("B26").Value = '3,42 '
$ws.Range("C26").Value = '3,42 '
$ws.Range("D26").Value = '1,13 '

# Row 27
$ws.Range("A27").Value = '20.02.2025'
$ws.Range("B27").Value = '3,43 '
$ws.Range("C27").Value = '3,43 '
$ws.Range("D27").Value = '1,13 '

# Row 28
$ws.Range("A28").Value = '21.02.2025'
$ws.Range("B28").Value = '3,42 '
$ws.Range("C28").Value = '3,42 '
$ws.Range("D28").Value = '1,13 '

# Row 29
$ws.Range("A29").Value = '22.02.2025'
$ws.Range("B29").Value = '3,42 '
$ws.Range("C29").Value = '3,42 '
$ws.Range("D29").Value = '1,13 '

# Row 30
$ws.Range("A30").Value = '23.02.2025'
$ws.Range("B30").Value = '3,42 '
$ws.Range("C30").Value = '3,42 '
$ws.Range("D30").Value = '1,13 '

# Row 31
$ws.Range("A31").Value = '24.02.2025'
$ws.Range("B31").Value = '3,42 '
$ws.Range("C31").Value = '3,42 '
$ws.Range("D31").Value = '1,13 '

# Row 32
$ws.Range("A32").Value = '25.02.2025'
$ws.Range("B32").Value = '3,42 '
$ws.Range("C32").Value = '3,42 '
$ws.Range("D32").Value = '1,13 '

